$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 208-209, shifting existing rows 208..297 down to 210..299.
$ws.Range("A208:R209").EntireRow.Insert()

# Populate new row 208 with the new record.
$ws.Range("A208").Value = 3
$ws.Range("B208").Value = "Femacal de La Calera"
$ws.Range("C208").Value = "Coquimbo"
$ws.Range("D208").Value = 44466
$ws.Range("E208").Value = 5
$ws.Range("F208").Value = 100114001
$ws.Range("G208").Value = "Papa"
$ws.Range("H208").Value = "Asterix"
$ws.Range("I208").Value = "1a (guarda)"
$ws.Range("J208").Value = 370
$ws.Range("K208").Value = 9500
$ws.Range("L208").Value = 10000
$ws.Range("M208").Value = 9757
$ws.Range("N208").Value = "$/saco 25 kilos"
$ws.Range("O208").Value = "Provincia de Talca"
$ws.Range("P208").Value = 390
$ws.Range("Q208").Value = 25
$ws.Range("R208").Value = "Hortaliza"

# Populate new row 209 with the new record.
$ws.Range("A209").Value = 3
$ws.Range("B209").Value = "Femacal de La Calera"
$ws.Range("C209").Value = "Coquimbo"
$ws.Range("D209").Value = 44466
$ws.Range("E209").Value = 5
$ws.Range("F209").Value = 100114001
$ws.Range("G209").Value = "Papa"
$ws.Range("H209").Value = "Rosara"
$ws.Range("I209").Value = "1a (guarda)"
$ws.Range("J209").Value = 340
$ws.Range("K209").Value = 10500
$ws.Range("L209").Value = 11000
$ws.Range("M209").Value = 10735
$ws.Range("N209").Value = "$/saco 25 kilos"
$ws.Range("O209").Value = "Provincia de Talca"
$ws.Range("P209").Value = 429
$ws.Range("Q209").Value = 25
$ws.Range("R209").Value = "Hortaliza"
